$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Week 5 entries (rows 17-19) ---
# Row 17
$ws.Range("A17").Value = 5
$ws.Range("B17").Value = 43697
$ws.Range("C17").Value = 0.39583333333333331
$ws.Range("D17").Value = 43697
$ws.Range("E17").Value = 0.4236111111111111
$ws.Range("F17").Value = "Formal team meeting"

# Row 18
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = 43699
$ws.Range("C18").Value = 0.875
$ws.Range("D18").Value = 43699
$ws.Range("E18").Value = 0.95833333333333337
$ws.Range("F18").Value = "Working on backend "

# Row 19
$ws.Range("A19").Value = 5
$ws.Range("B19").Value = 43700
$ws.Range("C19").Value = 0.83333333333333337
$ws.Range("D19").Value = 43700
$ws.Range("E19").Value = 0.875
$ws.Range("F19").Value = "Researching database integration"

# --- View: re-center zoom and move the active selection ---
$win = $ws.Application.ActiveWindow
$win.Zoom = 115
$ws.Range("F22").Select()
